$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.907.10"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "'1.668.13"
$ws.Range("E3").Value = "  +0.93%  "
$ws.Range("D5").Value = "'215.20"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").Value = "'0.520"
$ws.Range("E6").Value = "  +2.02%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  +1.32%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").Value = "'20.32"
$ws.Range("E10").Value = "  +1.91%  "
$ws.Range("E11").Value = "  +2.64%  "
$ws.Range("D12").Value = "'1.903.85"
$ws.Range("E12").Value = "  +0.92%  "
$ws.Range("D13").Value = "'1.690.37"
$ws.Range("E13").Value = "  +2.26%  "
$ws.Range("E14").Value = "  +0.05%  "
$ws.Range("E15").Value = "  +1.23%  "
$ws.Range("D16").Value = "'65.68"
$ws.Range("E16").Value = "  +0.52%  "
$ws.Range("D17").Value = "'26.905.50"
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("D18").Value = "'235.09"
$ws.Range("E18").Value = "  -2.05%  "
$ws.Range("D19").Value = "'7.96"
$ws.Range("E19").Value = "  +1.71%  "
$ws.Range("D20").Value = "'0.0₃0731"
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("D22").Value = "'4.39"
$ws.Range("E22").Value = "  -0.68%  "
$ws.Range("D23").Value = "'9.17"
$ws.Range("E23").Value = "  -0.95%  "
$ws.Range("E24").Value = "  -3.28%  "
$ws.Range("D25").Value = "'146.59"
$ws.Range("E25").Value = "  +0.33%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  -0.79%  "
$ws.Range("D28").Value = "'15.87"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  -0.70%  "
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("D33").Value = "'1.447.28"
$ws.Range("E33").Value = "  -4.83%  "
$ws.Range("E34").Value = "  +1.92%  "
$ws.Range("E35").Value = "  +2.71%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").Value = "'0.585"
$ws.Range("E37").Value = "  +0.86%  "
$ws.Range("D38").Value = "'0.903"
$ws.Range("E38").Value = "  +1.65%  "
$ws.Range("D39").Value = "'0.0170"
$ws.Range("E39").Value = "  +0.54%  "
$ws.Range("E40").Value = "  -3.81%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").Value = "'0.997"
$ws.Range("E42").Value = "  +8.97%  "
$ws.Range("E43").Value = "  +2.01%  "
$ws.Range("D44").Value = "'65.95"
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("D45").Value = "'1.809.25"
$ws.Range("E45").Value = "  +0.87%  "
$ws.Range("D46").Value = "'0.781"
$ws.Range("E46").Value = "  +0.99%  "
$ws.Range("D47").Value = "'90.78"
$ws.Range("E47").Value = "  +1.25%  "
$ws.Range("E48").Value = "  +0.91%  "
$ws.Range("E49").Value = "  -1.52%  "
$ws.Range("E50").Value = "  +4.11%  "
$ws.Range("E51").Value = "  -0.13%  "
